$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.726.80'
$ws.Range('E2').Value = '  -2.32%  '
$ws.Range('D3').Value = '1.875.65'
$ws.Range('E3').Value = '  -2.00%  '
$ws.Range('E4').Value = '  -0.87%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '247.47'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  +0.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.687'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  -2.99%  '
$ws.Range('E7').Value = '  -0.90%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.23'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').Value = '  +1.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.347'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('E9').Value = '  -2.64%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '50.71'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  -4.21%  '
$ws.Range('E11').Value = '  +0.75%  '
$ws.Range('E12').Value = '  -2.20%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '2.148.02'
$ws.Range('E13').Value = '  -1.97%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '12.83'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  +1.14%  '
$ws.Range('E15').Value = '  -0.69%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.89'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('D17').Value = '1.871.15'
$ws.Range('E17').Value = '  -2.29%  '
$ws.Range('D18').Value = '34.715.42'
$ws.Range('E18').Value = '  -2.30%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '72.90'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  -0.62%  '
$ws.Range('D20').Value = '0.0₃0824'
$ws.Range('E20').Value = '  -0.35%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '247.20'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  +1.84%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.73'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  -3.38%  '
$ws.Range('E23').Value = '  -3.40%  '
$ws.Range('E24').Value = '  -1.01%  '
$ws.Range('E25').Value = '  +3.67%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.24'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  -2.20%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '165.42'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  -1.97%  '
$ws.Range('E29').Value = '  -3.21%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.127'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  -3.97%  '
$ws.Range('D31').Value = '4.128.39'
$ws.Range('E31').Value = '  +0.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.67'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  +12.98%  '
$ws.Range('E33').Value = '  -0.39%  '
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.15'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  -1.63%  '
$ws.Range('E36').Value = '  -0.91%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.83'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  -5.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.835'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  -9.42%  '
$ws.Range('E39').Value = '  -3.76%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.38'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  -2.86%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '98.25'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  -0.54%  '
$ws.Range('E42').Value = '  +1.81%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0210'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  -0.19%  '
$ws.Range('E44').Value = '  -5.39%  '
$ws.Range('D45').Value = '1.292.91'
$ws.Range('E45').Value = '  -4.29%  '
$ws.Range('E46').Value = '  -4.46%  '
$ws.Range('E47').Value = '  -0.91%  '
$ws.Range('E48').Value = '  -2.17%  '
$ws.Range('E49').Value = '  +6.28%  '
$ws.Range('B50').Value = 'FraxShare'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.48'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  -1.40%  '
$ws.Range('B51').Value = 'Gas'
$ws.Range('C51').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '12.02'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  -1.14%  '
